$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells K1 / L1 -------------------------------------------------
$ws.Range("K1").Value = "fxppo2_accuracy_qkeras"
$ws.Range("L1").Value = "orig-fxppo2-drop_qkeras"

# Match the bold / bordered / centered header style used by the other header cells
foreach ($addr in @("K1", "L1")) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
    $cell.Borders.Item(8).LineStyle = 1   # xlEdgeTop
    $cell.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
    $cell.Borders.Item(10).LineStyle = 1  # xlEdgeRight
}

# --- New data columns K / L for rows 2..21 ------------------------------------
$kValues = @(
    0.8403755868544601,
    0.7777777777777778,
    0.7777777777777778,
    0.810641627543036,
    0.7777777777777778,
    0.8262910798122066,
    0.7777777777777778,
    0.7777777777777778,
    0.8075117370892019,
    0.7777777777777778,
    0.7777777777777778,
    0.7777777777777778,
    0.7777777777777778,
    0.8169014084507042,
    0.7777777777777778,
    0.7887323943661971,
    0.7777777777777778,
    0.7777777777777778,
    0.7777777777777778,
    0.8075117370892019
)

$lValues = @(
    0.00312989045383405,
    0,
    0,
    0.02660406885758992,
    0,
    0,
    0,
    0,
    0.00156494522691708,
    0,
    0,
    0.00469483568075113,
    0,
    0.00156494522691708,
    0,
    0.02816901408450712,
    -0.004694835680751241,
    0,
    0,
    -0.006259780907668211
)

for ($i = 0; $i -lt 20; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 11).Value = $kValues[$i]
    $ws.Cells.Item($row, 12).Value = $lValues[$i]
}
